$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Archive"
#
# The localization-status report is refreshed: the source file has moved
# from "Ready for handoff" into "In Translation", and its handback (target)
# has now actually landed, so each language sheet gains a populated
# "Latest Target File" / "Latest Handback File" pair plus a real
# "Latest Handback DateTime" (previously the zero-date placeholder).
# ---------------------------------------------------------------------------

$mdName  = "0194e2d8-cce1-4c92-b09a-100facd9bf7f.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/b7f2f4999d3b76643f9c91e74677b068478eda7f/e2e/0194e2d8-cce1-4c92-b09a-100facd9bf7f.md"

$zhName  = "0194e2d8-cce1-4c92-b09a-100facd9bf7f.6a513d531655ff2339c67e1066c71748c9d09f97.zh-cn.xlf"
$zhUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a2225a8a050856650829fd8bc4e1f77d3bd588b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0194e2d8-cce1-4c92-b09a-100facd9bf7f.6a513d531655ff2339c67e1066c71748c9d09f97.zh-cn.xlf"

$deName  = "0194e2d8-cce1-4c92-b09a-100facd9bf7f.6a513d531655ff2339c67e1066c71748c9d09f97.de-de.xlf"
$deUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/76a8d29621a87e92e233b33d5e1733754e8f9177/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0194e2d8-cce1-4c92-b09a-100facd9bf7f.6a513d531655ff2339c67e1066c71748c9d09f97.de-de.xlf"

# ---- Overview sheet ---------------------------------------------------------
# Mirrors the per-language Status column (B = zh-cn, C = de-de); both read the
# same "Ready for handoff" text as the language sheets, so it moves to
# "In Translation" in lock-step with them.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"

# ---- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "In Translation"

$wsZh.Range("E2").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, "", "", $mdName)

$wsZh.Range("F2").Value = $zhName
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhUrl, "", "", $zhName)

$wsZh.Range("G2").Value = "2016-02-23 08:17:22"
$wsZh.Range("H2").Value = "Include"

# ---- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "In Translation"

$wsDe.Range("E2").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, "", "", $mdName)

$wsDe.Range("F2").Value = $deName
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deUrl, "", "", $deName)

$wsDe.Range("G2").Value = "2016-02-23 08:17:43"
$wsDe.Range("H2").Value = "Include"
